# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cells (column Y:Z, row 1) used as a Text-formatted staging
# area so PasteSpecial(xlPasteValues) carries the literal digit string
# into the destination without touching the destination's existing
# (default) cell style -- avoids spuriously marking price/volume cells
# as explicit-General/quote-prefixed the way a direct .Value= would.
$scratch = $ws.Range("Y1:Z1")
$scratch.NumberFormat = "@"
$xlPasteValues = -4163

$ws.Range("Y1").Value = "26.838.89"
$ws.Range("Z1").Value = "  -1.98%  "
$scratch.Copy()
$ws.Range("D2:E2").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.809.44"
$ws.Range("Z1").Value = "  -1.26%  "
$scratch.Copy()
$ws.Range("D3:E3").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -0.35%  "
$ws.Range("Z1").Copy()
$ws.Range("E4").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "310.19"
$ws.Range("Z1").Value = "  -1.45%  "
$scratch.Copy()
$ws.Range("D5:E5").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -0.19%  "
$ws.Range("Z1").Copy()
$ws.Range("E6").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.4474"
$ws.Range("Z1").Value = "  +5.11%  "
$scratch.Copy()
$ws.Range("D7:E7").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.3661"
$ws.Range("Z1").Value = "  -1.18%  "
$scratch.Copy()
$ws.Range("D8:E8").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.07256"
$ws.Range("Z1").Value = "  -0.15%  "
$scratch.Copy()
$ws.Range("D9:E9").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.8500"
$ws.Range("Z1").Value = "  -2.07%  "
$scratch.Copy()
$ws.Range("D10:E10").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "20.64"
$ws.Range("Z1").Value = "  -2.50%  "
$scratch.Copy()
$ws.Range("D11:E11").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.810.26"
$ws.Range("Z1").Value = "  -1.28%  "
$scratch.Copy()
$ws.Range("D12:E12").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "6.599"
$ws.Range("Z1").Value = "  -2.11%  "
$scratch.Copy()
$ws.Range("D13:E13").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.07061"
$ws.Range("Z1").Value = "  -0.73%  "
$scratch.Copy()
$ws.Range("D14:E14").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "5.285"
$ws.Range("Z1").Value = "  -0.73%  "
$scratch.Copy()
$ws.Range("D15:E15").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "91.17"
$ws.Range("Z1").Value = "  +2.27%  "
$scratch.Copy()
$ws.Range("D16:E16").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.003"
$ws.Range("Z1").Value = "  -0.25%  "
$scratch.Copy()
$ws.Range("D17:E17").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.000008697"
$ws.Range("Z1").Value = "  -2.03%  "
$scratch.Copy()
$ws.Range("D18:E18").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.002"
$ws.Range("Z1").Value = "  -0.21%  "
$scratch.Copy()
$ws.Range("D19:E19").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "14.83"
$ws.Range("Z1").Value = "  -1.91%  "
$scratch.Copy()
$ws.Range("D20:E20").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "26.869.78"
$ws.Range("Z1").Value = "  -1.99%  "
$scratch.Copy()
$ws.Range("D21:E21").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "5.130"
$ws.Range("Z1").Value = "  -0.28%  "
$scratch.Copy()
$ws.Range("D22:E22").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "10.82"
$ws.Range("Z1").Value = "  -0.93%  "
$scratch.Copy()
$ws.Range("D23:E23").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.975"
$ws.Range("Z1").Value = "  -1.41%  "
$scratch.Copy()
$ws.Range("D24:E24").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "151.31"
$ws.Range("Z1").Value = "  -1.00%  "
$scratch.Copy()
$ws.Range("D25:E25").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "2.227"
$ws.Range("Z1").Value = "  +1.17%  "
$scratch.Copy()
$ws.Range("D26:E26").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "18.33"
$ws.Range("Z1").Value = "  -0.83%  "
$scratch.Copy()
$ws.Range("D27:E27").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "5.190"
$ws.Range("Z1").Value = "  -1.31%  "
$scratch.Copy()
$ws.Range("D28:E28").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "115.77"
$ws.Range("Z1").Value = "  -0.80%  "
$scratch.Copy()
$ws.Range("D29:E29").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.08807"
$ws.Range("Z1").Value = "  -0.99%  "
$scratch.Copy()
$ws.Range("D30:E30").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -3.14%  "
$ws.Range("Z1").Copy()
$ws.Range("E31").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.7444"
$ws.Range("Z1").Value = "  -2.42%  "
$scratch.Copy()
$ws.Range("D32:E32").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "2.932"
$ws.Range("Z1").Value = "  +4.04%  "
$scratch.Copy()
$ws.Range("D33:E33").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "4.418"
$ws.Range("Z1").Value = "  -1.70%  "
$scratch.Copy()
$ws.Range("D34:E34").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.002"
$ws.Range("Z1").Value = "  -0.28%  "
$scratch.Copy()
$ws.Range("D35:E35").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.087"
$ws.Range("Z1").Value = "  -3.30%  "
$scratch.Copy()
$ws.Range("D36:E36").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -1.75%  "
$ws.Range("Z1").Copy()
$ws.Range("E37").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.05180"
$ws.Range("Z1").Value = "  -2.28%  "
$scratch.Copy()
$ws.Range("D38:E38").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.5256"
$ws.Range("Z1").Value = "  +3.07%  "
$scratch.Copy()
$ws.Range("D39:E39").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "2.866"
$ws.Range("Z1").Value = "  -0.78%  "
$scratch.Copy()
$ws.Range("D40:E40").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "7.071"
$ws.Range("Z1").Value = "  -2.61%  "
$scratch.Copy()
$ws.Range("D41:E41").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.1683"
$ws.Range("Z1").Value = "  -1.57%  "
$scratch.Copy()
$ws.Range("D42:E42").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.5143"
$ws.Range("Z1").Value = "  +7.37%  "
$scratch.Copy()
$ws.Range("D43:E43").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "8.423"
$ws.Range("Z1").Value = "  -3.69%  "
$scratch.Copy()
$ws.Range("D44:E44").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "10.56"
$ws.Range("Z1").Value = "  -1.38%  "
$scratch.Copy()
$ws.Range("D45:E45").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "1.947"
$ws.Range("Z1").Value = "  +4.55%  "
$scratch.Copy()
$ws.Range("D46:E46").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "105.20"
$ws.Range("Z1").Value = "  -2.75%  "
$scratch.Copy()
$ws.Range("D47:E47").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -0.28%  "
$ws.Range("Z1").Copy()
$ws.Range("E48").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.06351"
$ws.Range("Z1").Value = "  -0.74%  "
$scratch.Copy()
$ws.Range("D49:E49").PasteSpecial($xlPasteValues)

$ws.Range("Z1").Value = "  -1.46%  "
$ws.Range("Z1").Copy()
$ws.Range("E50").PasteSpecial($xlPasteValues)

$ws.Range("Y1").Value = "0.9119"
$ws.Range("Z1").Value = "  -1.13%  "
$scratch.Copy()
$ws.Range("D51:E51").PasteSpecial($xlPasteValues)

# Clean up the scratch area
$excel.CutCopyMode = $false
$scratch.Clear()

Write-Host "Updated cryptos price/volume data."
